$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.69569456577301
$ws.Range("B1").Value = 4.663414001464844
$ws.Range("C1").Value = 2.914740562438965
$ws.Range("D1").Value = 1.490396976470947
$ws.Range("E1").Value = 1.09908664226532
